# Refresh the live crypto-price snapshot (Coin / Link / Price / Volume(1h))
# columns for rows 2..51, matching the upstream GitHub Actions data pull.
#
# NOTE: several Price values look numeric ("1.002", "0.05650", "236.51", ...)
# but must stay literal text (leading zeros / trailing zeros preserved, just
# like the source cells). A leading apostrophe is Excel's normal text-entry
# prefix, so it is prepended before assigning those values to .Value - exactly
# what typing them into Excel by hand would do.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.228.08"
$ws.Range("E2").Value = "  +0.24%  "
# Row 3
$ws.Range("D3").Value = "1.861.50"
$ws.Range("E3").Value = "  -0.10%  "
# Row 4
$ws.Range("E4").Value = "  +0.21%  "
# Row 5
$ws.Range("D5").Value = "'236.51"
$ws.Range("E5").Value = "  +1.07%  "
# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.15%  "
# Row 7
$ws.Range("D7").Value = "'0.4675"
$ws.Range("E7").Value = "  +0.16%  "
# Row 8
$ws.Range("D8").Value = "'0.2853"
$ws.Range("E8").Value = "  +0.66%  "
# Row 9
$ws.Range("D9").Value = "'0.06525"
$ws.Range("E9").Value = "  -0.41%  "
# Row 10
$ws.Range("D10").Value = "'21.97"
$ws.Range("E10").Value = "  +8.83%  "
# Row 11
$ws.Range("D11").Value = "'0.07915"
$ws.Range("E11").Value = "  +1.25%  "
# Row 12
$ws.Range("D12").Value = "'97.43"
$ws.Range("E12").Value = "  +1.32%  "
# Row 13
$ws.Range("D13").Value = "1.866.99"
$ws.Range("E13").Value = "  +0.63%  "
# Row 14
$ws.Range("D14").Value = "'5.160"
$ws.Range("E14").Value = "  +0.61%  "
# Row 15
$ws.Range("D15").Value = "'0.6810"
$ws.Range("E15").Value = "  +1.63%  "
# Row 16
$ws.Range("D16").Value = "'271.26"
$ws.Range("E16").Value = "  -3.13%  "
# Row 17
$ws.Range("D17").Value = "30.231.44"
$ws.Range("E17").Value = "  +0.18%  "
# Row 18
$ws.Range("D18").Value = "'13.53"
$ws.Range("E18").Value = "  +7.14%  "
# Row 19
$ws.Range("E19").Value = "  +0.05%  "
# Row 20
$ws.Range("D20").Value = "'0.000007345"
$ws.Range("E20").Value = "  +1.24%  "
# Row 21
$ws.Range("D21").Value = "2.112.85"
$ws.Range("E21").Value = "  +0.68%  "
# Row 22
$ws.Range("D22").Value = "'5.321"
$ws.Range("E22").Value = "  -2.46%  "
# Row 23
$ws.Range("E23").Value = "  +0.22%  "
# Row 24
$ws.Range("D24").Value = "'6.172"
$ws.Range("E24").Value = "  +0.46%  "
# Row 25
$ws.Range("D25").Value = "'167.76"
$ws.Range("E25").Value = "  +1.45%  "
# Row 26
$ws.Range("D26").Value = "'9.216"
$ws.Range("E26").Value = "  -1.02%  "
# Row 27
$ws.Range("D27").Value = "'18.94"
$ws.Range("E27").Value = "  +0.12%  "
# Row 28
$ws.Range("D28").Value = "'1.947"
$ws.Range("E28").Value = "  +2.26%  "
# Row 29
$ws.Range("D29").Value = "'1.385"
$ws.Range("E29").Value = "  +3.00%  "
# Row 30
$ws.Range("D30").Value = "'0.09811"
$ws.Range("E30").Value = "  +2.11%  "
# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.483"
$ws.Range("E31").Value = "  +0.88%  "
# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.356"
$ws.Range("E32").Value = "  -1.09%  "
# Row 33
$ws.Range("D33").Value = "'4.055"
$ws.Range("E33").Value = "  -1.55%  "
# Row 34
$ws.Range("D34").Value = "'0.04710"
$ws.Range("E34").Value = "  +1.00%  "
# Row 35
$ws.Range("D35").Value = "'1.132"
$ws.Range("E35").Value = "  +3.10%  "
# Row 36
$ws.Range("D36").Value = "'0.7015"
$ws.Range("E36").Value = "  +0.03%  "
# Row 37
$ws.Range("D37").Value = "'2.710"
$ws.Range("E37").Value = "  -0.01%  "
# Row 38
$ws.Range("D38").Value = "'0.01874"
$ws.Range("E38").Value = "  +0.67%  "
# Row 39
$ws.Range("D39").Value = "'2.629"
$ws.Range("E39").Value = "  +3.83%  "
# Row 40
$ws.Range("D40").Value = "'6.271"
$ws.Range("E40").Value = "  -0.09%  "
# Row 41
$ws.Range("D41").Value = "'75.30"
$ws.Range("E41").Value = "  +3.93%  "
# Row 42
$ws.Range("D42").Value = "'1.943"
$ws.Range("E42").Value = "  +0.84%  "
# Row 43
$ws.Range("D43").Value = "'0.8515"
$ws.Range("E43").Value = "  -0.06%  "
# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.07%  "
# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4157"
$ws.Range("E45").Value = "  -0.11%  "
# Row 46
$ws.Range("D46").Value = "'103.09"
$ws.Range("E46").Value = "  -0.06%  "
# Row 47
$ws.Range("D47").Value = "'7.173"
$ws.Range("E47").Value = "  +0.42%  "
# Row 48
$ws.Range("D48").Value = "'950.23"
$ws.Range("E48").Value = "  -3.83%  "
# Row 49
$ws.Range("D49").Value = "'9.250"
$ws.Range("E49").Value = "  +0.45%  "
# Row 50
$ws.Range("D50").Value = "'34.11"
$ws.Range("E50").Value = "  -0.15%  "
# Row 51
$ws.Range("D51").Value = "'0.05650"
$ws.Range("E51").Value = "  +0.19%  "
